$wb = $excel.ActiveWorkbook

# Add data for 2023-06-09: update 2023 (column J) violent crime counts
# across the citywide totals, by-neighborhood summary, and each neighborhood sheet.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 3133
$ws.Range("J3").Value = 3266
$ws.Range("J4").Value = 732
$ws.Range("J5").Value = 256
$ws.Range("J6").Value = 3865
$ws.Range("J7").Value = 11252

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 338
$ws.Range("J8").Value = 717
$ws.Range("J14").Value = 45
$ws.Range("J15").Value = 130
$ws.Range("J18").Value = 115
$ws.Range("J19").Value = 356
$ws.Range("J20").Value = 236
$ws.Range("J23").Value = 116
$ws.Range("J25").Value = 62
$ws.Range("J27").Value = 67
$ws.Range("J29").Value = 646
$ws.Range("J31").Value = 84
$ws.Range("J33").Value = 483
$ws.Range("J35").Value = 13
$ws.Range("J36").Value = 167
$ws.Range("J37").Value = 362
$ws.Range("J42").Value = 449
$ws.Range("J50").Value = 65
$ws.Range("J51").Value = 149
$ws.Range("J52").Value = 301
$ws.Range("J54").Value = 214
$ws.Range("J63").Value = 49
$ws.Range("J65").Value = 301
$ws.Range("J67").Value = 403
$ws.Range("J72").Value = 44
$ws.Range("J76").Value = 162
$ws.Range("J78").Value = 149
$ws.Range("J79").Value = 334
$ws.Range("J83").Value = 257
$ws.Range("J84").Value = 99
$ws.Range("J85").Value = 515
$ws.Range("J89").Value = 128
$ws.Range("J90").Value = 133
$ws.Range("J94").Value = 102
$ws.Range("J96").Value = 125
$ws.Range("J97").Value = 67
$ws.Range("J101").Value = 11252

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 194
$ws.Range("J4").Value = 39
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 515

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J5").Value = 24
$ws.Range("J6").Value = 212
$ws.Range("J7").Value = 717

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 113
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 44
$ws.Range("J4").Value = 13
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 17
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J5").Value = 15
$ws.Range("J7").Value = 362

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 168
$ws.Range("J7").Value = 403

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 74
$ws.Range("J3").Value = 97
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 257

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 128
$ws.Range("J3").Value = 153
$ws.Range("J6").Value = 159
$ws.Range("J7").Value = 483

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 54
$ws.Range("J7").Value = 214

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 194
$ws.Range("J3").Value = 222
$ws.Range("J4").Value = 38
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 646

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 100
$ws.Range("J6").Value = 139
$ws.Range("J7").Value = 356

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 92
$ws.Range("J3").Value = 98
$ws.Range("J7").Value = 449

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 50
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 98
$ws.Range("J3").Value = 121
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 73
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 236

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 34
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 115

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 60
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 102

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 40
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 13

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 16
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 43
$ws.Range("J4").Value = 4
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J4").Value = 19
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 44
